# The template paragraph holds a single Word field whose code is the M2Doc
# query "m:Sequence{...}" written out as fldChar(begin)/instrText*/fldChar(end).
# The commit switches the parser to a token-iterator based rewriter that
# expects the query spelled out as plain run text (a "{...}" block) instead
# of a real field, so we replace the field's run sequence with equivalent
# <w:t> runs carrying the same text (plus the leading "{" and trailing "}"
# that used to be implicit in the field delimiters), while keeping the
# _GoBack bookmark exactly where it was.

$d = $word.ActiveDocument

# Locate the field and the paragraph that fully contains it.
$field = $d.Fields(1)
$codeStart = $field.Code.Start
$target = $null
foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    if (($codeStart -ge $pr.Start) -and ($codeStart -lt $pr.End)) {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the paragraph containing the M2Doc field"
}

$r = $target.Range

# Build the replacement run sequence: every former instrText chunk becomes a
# <w:t> run with identical text, the opening fldChar+" " become a single "{"
# run, the closing " "+fldChar become a single "}" run, and the _GoBack
# bookmark is preserved in its original position (between "even more" and
# " text'.sampleText(").
$innerXml = '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>Sequence{</w:t></w:r>' +
    '<w:r><w:t>''some text''.sampleText(6)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:r><w:t>''</w:t></w:r>' +
    '<w:r><w:t>more</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> text''.sampleText(</w:t></w:r>' +
    '<w:r><w:t>8</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:r><w:t>''</w:t></w:r>' +
    '<w:r><w:t>even more</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve"> text''.sampleText(</w:t></w:r>' +
    '<w:r><w:t>10</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '<w:r><w:t>}</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

Write-Host "Field converted to plain text runs."
